# Committing final code for Job description page
# Adds a second username/password pair (row 3) with mailto hyperlinks,
# matching the existing row-2 convention, plus a trailing blank row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New credential values for row 3
$ws.Range("A3").Value = "saamaajik.fed.cci@gmail.com"
$ws.Range("B3").Value = "S0c!al@2022"

# Hyperlink each new cell to its own value via a mailto: link, same as A2/B2
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:saamaajik.fed.cci@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:S0c!al@2022")

# Copy the look (border + hyperlink font/style) from row 2 onto row 3
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Add a trailing blank bordered row 4, matching the look of row 1 (no hyperlink font)
$ws.Range("A1:B1").Copy()
$ws.Range("A4:B4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("A4:B4").Value = ""

# Restore the user's on-screen selection to J17
$ws.Range("J17").Select()
